$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Flow")

# Fill in the "Seq" numbers (column A) for the existing rows 15-20
$ws.Cells.Item(15, 1).Value = 31
$ws.Cells.Item(16, 1).Value = 32
$ws.Cells.Item(17, 1).Value = 33
$ws.Cells.Item(18, 1).Value = 34
$ws.Cells.Item(19, 1).Value = 35
$ws.Cells.Item(20, 1).Value = 36

# Add the new row describing the catch code in the Data/RefreshTodaysMatchups method
$ws.Cells.Item(21, 1).Value = 40
$ws.Cells.Item(21, 2).Value = "Select Game Date"
$ws.Cells.Item(21, 3).Value = "RefreshTodaysMatchups()"

# Match the look of the existing "RefreshTodaysMatchups()" code cell (C5) by
# copying its formatting (font, etc.) onto the new cell
$ws.Cells.Item(5, 3).Copy()
$ws.Cells.Item(21, 3).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the selection / active cell shown in the saved workbook
$ws.Range("G9").Select()
